$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.273.53'
$ws.Range('E2').Value = '  -2.90%  '
$ws.Range('D3').Value = '2.228.79'
$ws.Range('E3').Value = '  -4.65%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.79'
$ws.Range('E5').Value = '  -3.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.84'
$ws.Range('E6').Value = '  -6.14%  '
$ws.Range('E7').Value = '  -7.50%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -7.44%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.11'
$ws.Range('E10').Value = '  -7.59%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.11'
$ws.Range('E11').Value = '  -2.34%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0829'
$ws.Range('E12').Value = '  -8.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.70'
$ws.Range('E13').Value = '  -7.13%  '
$ws.Range('E14').Value = '  -1.51%  '
$ws.Range('D15').Value = '2.568.63'
$ws.Range('E15').Value = '  -4.61%  '
$ws.Range('E16').Value = '  -10.40%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.37'
$ws.Range('E17').Value = '  -5.39%  '
$ws.Range('D18').Value = '2.229.55'
$ws.Range('E18').Value = '  -4.48%  '
$ws.Range('D19').Value = '43.200.97'
$ws.Range('E19').Value = '  -3.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.89'
$ws.Range('E20').Value = '  -9.04%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '0.0₃0968'
$ws.Range('E21').Value = '  -8.04%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.56'
$ws.Range('E22').Value = '  -9.03%  '
$ws.Range('E23').Value = '  -10.77%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.19'
$ws.Range('E24').Value = '  -10.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '236.89'
$ws.Range('E25').Value = '  -7.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.19'
$ws.Range('E26').Value = '  -2.64%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '4.06'
$ws.Range('E28').Value = '  +1.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.04'
$ws.Range('E29').Value = '  -9.84%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.43'
$ws.Range('E31').Value = '  -11.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '36.62'
$ws.Range('E32').Value = '  -0.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.25'
$ws.Range('E33').Value = '  -7.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0864'
$ws.Range('E34').Value = '  -7.98%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '157.00'
$ws.Range('E35').Value = '  -4.87%  '
$ws.Range('E36').Value = '  -3.91%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.29'
$ws.Range('E37').Value = '  +3.01%  '
$ws.Range('E38').Value = '  -7.44%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.42'
$ws.Range('E39').Value = '  -5.09%  '
$ws.Range('B40').Value = 'ARBITRUM'
$ws.Range('C40').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.85'
$ws.Range('E40').Value = '  -3.32%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.105'
$ws.Range('E41').Value = '  -9.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.70'
$ws.Range('E42').Value = '  -4.82%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0318'
$ws.Range('E43').Value = '  -8.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.25'
$ws.Range('E44').Value = '  +11.74%  '
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('D46').Value = '1.769.03'
$ws.Range('E46').Value = '  -5.34%  '
$ws.Range('E47').Value = '  -9.65%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '83.71'
$ws.Range('E48').Value = '  -11.47%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.28'
$ws.Range('E49').Value = '  -11.81%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.86'
$ws.Range('E50').Value = '  -3.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.83'
$ws.Range('E51').Value = '  -10.84%  '
